# Finalize production module KL-E: add two new worksheets (elasKL-E, prodKL-E)
# after prodKL, mirroring elasKL / prodKL but with an extra "energy" column,
# and move the selected/active tab to the new last sheet (prodKL-E).

$wb = $excel.ActiveWorkbook

# Country / region codes shared by every data row (rows 2-36) on both
# new sheets - same order as on the existing elasKL / prodKL sheets.
$codes = @(
    "iPARI","iWHEA","iOCER","iFVEG","iOILS","iSUGB","iFIBR","iOTHC","iANIM",
    "iFORE","iFISH","iFOSM","iOTHM","iFBTO","iTXWO","iCOKE","iREFN","iCHEM",
    "iRUBP","iNMMP","iMETP","iELEC","iMACH","iELCF","iELCG","iTRDI","iHWAT",
    "iWATR","iCONS","iTRAD","iHORE","iTRAN","iREBA","iPUBO","iWAST"
)

# --- New sheet 1: elasKL-E (inserted right after prodKL) -------------------
$wsElas = $wb.Worksheets.Add($null, $wb.Worksheets("prodKL"))
$wsElas.Name = "elasKL-E"

$wsElas.Range("B1").Value = "elasKLE"
$wsElas.Range("C1").Value = "elasKL"
$wsElas.Range("D1").Value = "elasE"

for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 2
    $wsElas.Cells.Item($r, 1).Value = "'" + $codes[$i]
    $wsElas.Cells.Item($r, 2).Value = 0.4
    $wsElas.Cells.Item($r, 3).Value = 0.95
    $wsElas.Cells.Item($r, 4).Value = 0.5
}

# --- New sheet 2: prodKL-E (inserted right after elasKL-E, becomes active) -
$wsProd = $wb.Worksheets.Add($null, $wsElas)
$wsProd.Name = "prodKL-E"

$wsProd.Range("B1").Value = "COE"
$wsProd.Range("C1").Value = "GOS"
$wsProd.Range("D1").Value = "ENER"

for ($i = 0; $i -lt $codes.Length; $i++) {
    $r = $i + 2
    $wsProd.Cells.Item($r, 1).Value = "'" + $codes[$i]
    $wsProd.Cells.Item($r, 2).Value = 1
    $wsProd.Cells.Item($r, 3).Value = 1
    $wsProd.Cells.Item($r, 4).Value = 1
}

# prodKL-E is the new last / active tab.
$wsProd.Select()
